# eventannotations.xlsx edit:
#  - Insert 3 new rows (Description/Type/SliderUnit fields) into the
#    "Annotations" sheet right before the "~dCHARGE" row.
#  - Make "Annotations" the active/selected sheet (it was "Examples").
#  - Reset selections on both sheets to match the new layout.

$wb = $excel.ActiveWorkbook
$annotations = $wb.Worksheets.Item("Annotations")
$examples    = $wb.Worksheets.Item("Examples")

# Insert three blank rows above the old row 6 ("~dCHARGE ...").
$annotations.Rows("6:8").Insert()

# New row 6: ~descr
$annotations.Range("A6").Value = "~descr"
$annotations.Range("B6").Value = "The Description field of the Event"
$annotations.Range("C6").Value = "Gas 10"

# New row 7: ~type
$annotations.Range("A7").Value = "~type"
$annotations.Range("B7").Value = "The Type field of the Event"
$annotations.Range("C7").Value = "Power"

# New row 8: ~sldrunit
$annotations.Range("A8").Value = "~sldrunit"
$annotations.Range("B8").Value = "The value of the Slider Unit for this Event"
$annotations.Range("C8").Value = "kPa"

# Examples sheet keeps its own cursor at A9 but is no longer the active tab.
$examples.Range("A9").Select()

# Make Annotations the active sheet (was Examples) and move its selection to A9.
$annotations.Activate()
$annotations.Range("A9").Select()
